# GNUPlot e Gráficos ISCC
#
# Update the multiplier constants used in the "proportion of coincidents"
# (proporção de coincidentes) calculations. Each block multiplies a
# different sample size by the ratio of the "new" occurrence count to the
# baseline total; the multiplier is being corrected per block.
#
# I3  (block "G2"):  306 -> 272
# O3  (block "M2"):  306 -> 182
# D16 (block "B15"): 306 -> 132
# I23 (block "G22"): 306 -> 90
#
# The dependent summary cells (I4, O4, D17, I24) recompute automatically
# because they reference the cells above via formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Formula = "=272*(H3/G2)"
$ws.Range("O3").Formula = "=182*(N3/M2)"
$ws.Range("D16").Formula = "=132*(C16/B15)"
$ws.Range("I23").Formula = "=90*(H23/G22)"

# Move the viewport/selection to match where the author ended up working.
$excel.Goto($ws.Range("A7"), $true)
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("I24").Select()
